$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Task 16 note: append "Task 18: Complete" note
$ws.Range("F14").Value = "Task 16: Complete. Task 18: Complete"

# E14 hours: 2 -> 4
$ws.Range("E14").Value = 4

# Update selection to F15
$ws.Range("F15").Select()
